# Updated R code and solution based on given constraint
$wb = $excel.ActiveWorkbook

# --- "Table - 2.1" sheet: the suitable investment type given the funding
# constraint changes from the stray "post_ipo_debt?" placeholder to the
# correct answer "venture".
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C9").Value = "venture"
$null = $ws2.Range("H9").Select()

# --- "Table-3.1" sheet: second/third top English speaking country answers
# are corrected from Canada/New Zealand to United Kingdom/India (the USA
# answer in C5 stays the same, only its shared-string index shifts because
# the now-unused "post_ipo_debt?" string was dropped from the table).
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C6").Value = "United Kingdom of Great Britain and Northern Ireland"
$ws3.Range("C7").Value = "India"

# Column C widens to fit the longer "United Kingdom ..." answer.
$ws3.Columns.Item(3).ColumnWidth = 40.33

$null = $ws3.Range("C8").Select()
